$newValues = @(
    "9+64=",
    "4+87=",
    "27-9=",
    "62-18=",
    "90-89=",
    "38-29=",
    "68+7=",
    "17+16=",
    "54+18=",
    "35+8=",
    "90-13=",
    "38+6=",
    "83+8=",
    "73-38=",
    "43-35=",
    "39+26=",
    "31-9=",
    "46+8=",
    "46-38=",
    "26+39=",
    "36+26=",
    "57+19=",
    "35+46=",
    "94-49=",
    "44+48=",
    "59+24=",
    "92-45=",
    "6+46=",
    "57+7=",
    "6+67=",
    "78+5=",
    "12-4=",
    "70-61=",
    "14+79=",
    "91-38=",
    "77-69=",
    "51-13=",
    "32-16=",
    "91-59=",
    "27+55=",
    "60-5=",
    "16+57=",
    "7+19=",
    "23-9=",
    "69+24=",
    "88-39=",
    "96-68=",
    "47+24=",
    "45+39=",
    "69+18=",
    "12+69=",
    "27+37=",
    "18+4=",
    "77+17=",
    "76-7=",
    "63+8=",
    "52-24=",
    "18+43=",
    "56-49=",
    "88+3=",
    "91-64=",
    "18+8=",
    "55-27=",
    "24-19=",
    "91-87=",
    "83-74=",
    "47+39=",
    "8+16=",
    "9+64=",
    "46-7=",
    "16+9=",
    "19+64=",
    "17+45=",
    "28+8=",
    "4+57=",
    "33-8=",
    "73-6=",
    "84-7=",
    "38+28=",
    "66+15=",
    "82-39=",
    "48+24=",
    "29+52=",
    "82-19=",
    "16+28=",
    "73-17=",
    "16+59=",
    "56-38=",
    "29+38=",
    "19+7=",
    "28+23=",
    "26+18=",
    "16+68=",
    "47+49=",
    "51-13=",
    "45+39=",
    "13-4=",
    "19+68=",
    "70-8=",
    "31-13="
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells"
